$d = $word.ActiveDocument

# Locate the target paragraph by its current text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Realizar abertura de uma m*quina.*") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Target paragraph not found"
}

$r = $target.Range

$newParaXml = '<w:p w:rsidR="00D36C55" w:rsidRPr="007E4E0C" w:rsidRDefault="00D36C55" w:rsidP="00B55D22"><w:pPr><w:pStyle w:val="Cabealho"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="007E4E0C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">REALIZAR ABERTURA DE </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="007E4E0C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>UMA MÁQUINA.</w:t></w:r></w:p>'

$pkgXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($pkgXml)

Write-Output "After: $($target.Range.Text)"
